# Auto commit refresh: new metrics data lands on the "Metrics" sheet; the
# "today" sheet just re-derives its B/E/F columns from Metrics via existing
# formulas, so writing the source values and letting the workbook recalc
# is enough to keep everything consistent.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# --- refresh the source figures on Metrics!B2:B13 -------------------------
$metrics.Range("B2").Value  = 392469.8000000001
$metrics.Range("B3").Value  = 345845.22000000003
$metrics.Range("B4").Value  = 120783.14
$metrics.Range("B5").Value  = 15949
$metrics.Range("B6").Value  = 5188715.5500000007
$metrics.Range("B7").Value  = 4387921.9000000013
$metrics.Range("B8").Value  = 1527742.9700000004
$metrics.Range("B9").Value  = 202156
$metrics.Range("B10").Value = 33654096.540000014
$metrics.Range("B11").Value = 31663197.060000006
$metrics.Range("B12").Value = 11809465.009999998
$metrics.Range("B13").Value = 1299786

# today!B11:B22 / E11:E22 / F11:F22 are plain formulas that reference the
# Metrics cells above (and A1 is a volatile TODAY()-1), so recalculating
# the workbook brings them back in sync without touching them directly.
$excel.Calculate()

# --- restore the cursor/selection state recorded in the saved view -------
[void]$metrics.Activate()
[void]$metrics.Range("E21").Select()

[void]$today.Activate()
[void]$today.Range("I7").Select()
